$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 'ECs'
$ws.Cells.Item(2, 2).Value = 'Col1a1'
$ws.Cells.Item(2, 3).Value = 'Gp6'
$ws.Cells.Item(2, 4).Value = 'FAPs'
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 91.94136433333334
$ws.Cells.Item(2, 8).Value = 275.824093
$ws.Cells.Item(2, 9).Value = 0.02307547609860541
$ws.Cells.Item(2, 10).Value = 0.02307547609860541
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.009511333333333333
$ws.Cells.Item(2, 14).Value = 0.028534
$ws.Cells.Item(2, 15).Value = 0.7967497835981349
$ws.Cells.Item(2, 16).Value = 0.7967497835981349
$ws.Cells.Item(2, 17).Value = 0.8744849632957779
$ws.Cells.Item(2, 18).Value = 7.870364669662
$ws.Cells.Item(2, 19).Value = 0.0183853805879878
$ws.Cells.Item(2, 20).Value = 0.0183853805879878

# Row 3
$ws.Cells.Item(3, 1).Value = 'ECs'
$ws.Cells.Item(3, 2).Value = 'Col1a1'
$ws.Cells.Item(3, 3).Value = 'Gp6'
$ws.Cells.Item(3, 4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 91.94136433333334
$ws.Cells.Item(3, 8).Value = 275.824093
$ws.Cells.Item(3, 9).Value = 0.02307547609860541
$ws.Cells.Item(3, 10).Value = 0.02307547609860541
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.002426333333333333
$ws.Cells.Item(3, 14).Value = 0.007279
$ws.Cells.Item(3, 15).Value = 0.2032502164018652
$ws.Cells.Item(3, 16).Value = 0.2032502164018652
$ws.Cells.Item(3, 17).Value = 0.2230803969941111
$ws.Cells.Item(3, 18).Value = 2.007723572947
$ws.Cells.Item(3, 19).Value = 0.004690095510617619
$ws.Cells.Item(3, 20).Value = 0.004690095510617619

# Row 4
$ws.Cells.Item(4, 1).Value = 'FAPs'
$ws.Cells.Item(4, 2).Value = 'Col1a1'
$ws.Cells.Item(4, 3).Value = 'Gp6'
$ws.Cells.Item(4, 4).Value = 'FAPs'
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 3857.568359333333
$ws.Cells.Item(4, 8).Value = 11572.705078
$ws.Cells.Item(4, 9).Value = 0.9681738695089209
$ws.Cells.Item(4, 10).Value = 0.9681738695089209
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.009511333333333333
$ws.Cells.Item(4, 14).Value = 0.028534
$ws.Cells.Item(4, 15).Value = 0.7967497835981349
$ws.Cells.Item(4, 16).Value = 0.7967497835981349
$ws.Cells.Item(4, 17).Value = 36.6906185217391
$ws.Cells.Item(4, 18).Value = 330.215566695652
$ws.Cells.Item(4, 19).Value = 0.7713923210166016
$ws.Cells.Item(4, 20).Value = 0.7713923210166016

# Row 5
$ws.Cells.Item(5, 1).Value = 'FAPs'
$ws.Cells.Item(5, 2).Value = 'Col1a1'
$ws.Cells.Item(5, 3).Value = 'Gp6'
$ws.Cells.Item(5, 4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 3857.568359333333
$ws.Cells.Item(5, 8).Value = 11572.705078
$ws.Cells.Item(5, 9).Value = 0.9681738695089209
$ws.Cells.Item(5, 10).Value = 0.9681738695089209
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.002426333333333333
$ws.Cells.Item(5, 14).Value = 0.007279
$ws.Cells.Item(5, 15).Value = 0.2032502164018652
$ws.Cells.Item(5, 16).Value = 0.2032502164018652
$ws.Cells.Item(5, 17).Value = 9.359746695862443
$ws.Cells.Item(5, 18).Value = 84.237720262762
$ws.Cells.Item(5, 19).Value = 0.1967815484923194
$ws.Cells.Item(5, 20).Value = 0.1967815484923194

# Row 6
$ws.Cells.Item(6, 1).Value = 'Inflammatory-Mac'
$ws.Cells.Item(6, 2).Value = 'Col1a1'
$ws.Cells.Item(6, 3).Value = 'Gp6'
$ws.Cells.Item(6, 4).Value = 'FAPs'
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.840730666666667
$ws.Cells.Item(6, 8).Value = 5.522192
$ws.Cells.Item(6, 9).Value = 0.000461987232956876
$ws.Cells.Item(6, 10).Value = 0.000461987232956876
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.009511333333333333
$ws.Cells.Item(6, 14).Value = 0.028534
$ws.Cells.Item(6, 15).Value = 0.7967497835981349
$ws.Cells.Item(6, 16).Value = 0.7967497835981349
$ws.Cells.Item(6, 17).Value = 0.01750780294755556
$ws.Cells.Item(6, 18).Value = 0.157570226528
$ws.Cells.Item(6, 19).Value = 0.0003680882278834921
$ws.Cells.Item(6, 20).Value = 0.0003680882278834921

# Row 7
$ws.Cells.Item(7, 1).Value = 'Inflammatory-Mac'
$ws.Cells.Item(7, 2).Value = 'Col1a1'
$ws.Cells.Item(7, 3).Value = 'Gp6'
$ws.Cells.Item(7, 4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.840730666666667
$ws.Cells.Item(7, 8).Value = 5.522192
$ws.Cells.Item(7, 9).Value = 0.000461987232956876
$ws.Cells.Item(7, 10).Value = 0.000461987232956876
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(7, 12).Value = 0.3333333333333333
$ws.Cells.Item(7, 13).Value = 0.002426333333333333
$ws.Cells.Item(7, 14).Value = 0.007279
$ws.Cells.Item(7, 15).Value = 0.2032502164018652
$ws.Cells.Item(7, 16).Value = 0.2032502164018652
$ws.Cells.Item(7, 17).Value = 0.004466226174222222
$ws.Cells.Item(7, 18).Value = 0.040196035568
$ws.Cells.Item(7, 19).Value = 0.00009389900507338398
$ws.Cells.Item(7, 20).Value = 0.00009389900507338399

# Row 8
$ws.Cells.Item(8, 1).Value = 'MuSCs'
$ws.Cells.Item(8, 2).Value = 'Col1a1'
$ws.Cells.Item(8, 3).Value = 'Gp6'
$ws.Cells.Item(8, 4).Value = 'FAPs'
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 31.11921133333334
$ws.Cells.Item(8, 8).Value = 93.357634
$ws.Cells.Item(8, 9).Value = 0.007810310653280575
$ws.Cells.Item(8, 10).Value = 0.007810310653280575
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.009511333333333333
$ws.Cells.Item(8, 14).Value = 0.028534
$ws.Cells.Item(8, 15).Value = 0.7967497835981349
$ws.Cells.Item(8, 16).Value = 0.7967497835981349
$ws.Cells.Item(8, 17).Value = 0.2959851920617778
$ws.Cells.Item(8, 18).Value = 2.663866728556
$ws.Cells.Item(8, 19).Value = 0.006222863322835506
$ws.Cells.Item(8, 20).Value = 0.006222863322835506

# Row 9
$ws.Cells.Item(9, 1).Value = 'MuSCs'
$ws.Cells.Item(9, 2).Value = 'Col1a1'
$ws.Cells.Item(9, 3).Value = 'Gp6'
$ws.Cells.Item(9, 4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 31.11921133333334
$ws.Cells.Item(9, 8).Value = 93.357634
$ws.Cells.Item(9, 9).Value = 0.007810310653280575
$ws.Cells.Item(9, 10).Value = 0.007810310653280575
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.002426333333333333
$ws.Cells.Item(9, 14).Value = 0.007279
$ws.Cells.Item(9, 15).Value = 0.2032502164018652
$ws.Cells.Item(9, 16).Value = 0.2032502164018652
$ws.Cells.Item(9, 17).Value = 0.07550557976511112
$ws.Cells.Item(9, 18).Value = 0.679550217886
$ws.Cells.Item(9, 19).Value = 0.00158744733044507
$ws.Cells.Item(9, 20).Value = 0.00158744733044507

# Row 10
$ws.Cells.Item(10, 1).Value = 'Resolving-Mac'
$ws.Cells.Item(10, 2).Value = 'Col1a1'
$ws.Cells.Item(10, 3).Value = 'Gp6'
$ws.Cells.Item(10, 4).Value = 'FAPs'
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 1.905952
$ws.Cells.Item(10, 8).Value = 5.717856
$ws.Cells.Item(10, 9).Value = 0.0004783565062362683
$ws.Cells.Item(10, 10).Value = 0.0004783565062362683
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.009511333333333333
$ws.Cells.Item(10, 14).Value = 0.028534
$ws.Cells.Item(10, 15).Value = 0.7967497835981349
$ws.Cells.Item(10, 16).Value = 0.7967497835981349
$ws.Cells.Item(10, 17).Value = 0.01812814478933333
$ws.Cells.Item(10, 18).Value = 0.163153303104
$ws.Cells.Item(10, 19).Value = 0.0003811304428265067
$ws.Cells.Item(10, 20).Value = 0.0003811304428265067

# Row 11
$ws.Cells.Item(11, 1).Value = 'Resolving-Mac'
$ws.Cells.Item(11, 2).Value = 'Col1a1'
$ws.Cells.Item(11, 3).Value = 'Gp6'
$ws.Cells.Item(11, 4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 1.905952
$ws.Cells.Item(11, 8).Value = 5.717856
$ws.Cells.Item(11, 9).Value = 0.0004783565062362683
$ws.Cells.Item(11, 10).Value = 0.0004783565062362683
$ws.Cells.Item(11, 11).Value = 1
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.002426333333333333
$ws.Cells.Item(11, 14).Value = 0.007279
$ws.Cells.Item(11, 15).Value = 0.2032502164018652
$ws.Cells.Item(11, 16).Value = 0.2032502164018652
$ws.Cells.Item(11, 17).Value = 0.004624474869333334
$ws.Cells.Item(11, 18).Value = 0.041620273824
$ws.Cells.Item(11, 19).Value = 0.00009722606340976174
$ws.Cells.Item(11, 20).Value = 0.00009722606340976174
